$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data for column A (file names) and column C (extension) / D (folder) ---
$names = @(
    "[Daman] Elfen Lied 01 [1280x720_Blu_Ray_Dual_Audio_FLAC][71978425]",
    "[Daman] Elfen Lied 02 [1280x720_Blu_Ray_Dual_Audio_FLAC][1968c90e]",
    "[Daman] Elfen Lied 03 [1280x720_Blu_Ray_Dual_Audio_FLAC][85ac8d66]",
    "[Daman] Elfen Lied 04 [1280x720_Blu_Ray_Dual_Audio_FLAC][59abc835]",
    "[Daman] Elfen Lied 05 [1280x720_Blu_Ray_Dual_Audio_FLAC][33a67a3b]",
    "[Daman] Elfen Lied 06 [1280x720_Blu_Ray_Dual_Audio_FLAC][3ef05c2a]",
    "[Daman] Elfen Lied 07 [1280x720_Blu_Ray_Dual_Audio_FLAC][5d4c6a5a]",
    "[Daman] Elfen Lied 08 [1280x720_Blu_Ray_Dual_Audio_FLAC][a4abbe62]",
    "[Daman] Elfen Lied 09 [1280x720_Blu_Ray_Dual_Audio_FLAC][631d7313]",
    "[Daman] Elfen Lied 10 [1280x720_Blu_Ray_Dual_Audio_FLAC][c48df422]",
    "[Daman] Elfen Lied 11 [1280x720_Blu_Ray_Dual_Audio_FLAC][ba639ab7]",
    "[Daman] Elfen Lied 12 [1280x720_Blu_Ray_Dual_Audio_FLAC][e0889607]",
    "[Daman] Elfen Lied 13 [1280x720_Blu_Ray_Dual_Audio_FLAC][740bded2]",
    "[Daman] Elfen Lied OVA [1280x720_Blu_Ray_FLAC][e82ad193]"
)

$folder = "E:/TV&Movies/Anime/[Daman] Elfen Lied"

# Clear the old contents first (old sheet had rows 1-14, new has 1-15)
$ws.Cells.Clear()

# --- Header row (A, B, C already existed; D introduces the new "Folder" string) ---
$ws.Range("A1").Value = "File Name"
$ws.Range("B1").Value = "New File Name"
$ws.Range("C1").Value = "Extension"
$ws.Range("D1").Value = "Folder"

# --- Data rows 2..15 ---
# Write column A, then C, then D for each row (B is filled afterwards via formula)
# so that new shared-string entries are created in the same row-major order
# as the target workbook.
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = ".mkv"
    $ws.Cells.Item($row, 4).Value = $folder
}

# Column B: formulas that take the first 21 characters of column A.
# B2 is entered on its own, then B3:B15 is entered as one range so it is
# stored as a single shared-formula group (matches how this was authored).
$ws.Range("B2").Formula = "=LEFT(A2,21)"
$ws.Range("B3:B15").Formula = "=LEFT(A3,21)"

# --- View / column sizing ---
$ws.Range("B2:B15").Select()
$ws.Columns.Item(2).AutoFit() | Out-Null
